$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - B14's stored type changes from text "4" to a genuine number 4.
$ws.Range("B14").Value = 4

# Row 15 - a brand-new annotation row appended below the existing data.
$ws.Range("A15").Value = "Sunsi Wu"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "4"
$ws.Range("C15").Value = "well carried; very thorough"
$ws.Range("D15").Value = "APC"
$ws.Range("E15").Value = "EXP"
$ws.Range("F15").Value = "ea04c829-c996-4167-8585-03efb193cd41"
$ws.Range("G15").Value = "ByOExmWAb_annotated.xlsx"
$ws.Range("H15").Value = "The experiments were well carried through and very thorough."
